# Generate Report for Handoff
# Update status text from "In Translation" to "Ready for handoff" and
# refresh the handoff timestamps on the Overview / zh-cn / de-de sheets.
# Also widen the "Status" / language-status columns to fit the new text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-25 12:41:04"

# --- zh-cn sheet ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-25 12:40:56"

# --- de-de sheet ---
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-25 12:41:04"

# --- Widen the Status columns so the longer text fits ---
$wsOverview.Range("E1:F1").ColumnWidth = 17.2159881591797
$wsZhCn.Range("C1").ColumnWidth = 17.2159881591797
$wsDeDe.Range("C1").ColumnWidth = 17.2159881591797
